# Fix ACO Insights basic calc (do this because it was fixed by column).
$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsCalcs = $wb.Worksheets.Item("Calcs")

# Swap the ACO Insight toggle values in column E (Data Feed #2) between rows 18 and 19.
$wsInputs.Range("E18").Value = "Y"
$wsInputs.Range("E19").Value = "N"

# The row-61 formulas in Calcs were hard-coded to reference Inputs!$D$19 regardless
# of column; fix them to reference the row-19 cell in their own column instead.
$wsCalcs.Range("D61").Formula = '=IF(Inputs!D$19 = "Y", 0, IF(Inputs!D18="Y", Prices!$C20, 0))'
$wsCalcs.Range("E61").Formula = '=IF(Inputs!E$19 = "Y", 0, IF(Inputs!E18="Y", Prices!$C20, 0))'
$wsCalcs.Range("F61").Formula = '=IF(Inputs!F$19 = "Y", 0, IF(Inputs!F18="Y", Prices!$C20, 0))'
$wsCalcs.Range("G61").Formula = '=IF(Inputs!G$19 = "Y", 0, IF(Inputs!G18="Y", Prices!$C20, 0))'
$wsCalcs.Range("H61").Formula = '=IF(Inputs!H$19 = "Y", 0, IF(Inputs!H18="Y", Prices!$C20, 0))'

$excel.Calculate()
